$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "salutation [Link] [Salutation]"
$ws.Range("D1").Value = "customer_group [Link] [Customer Group]"
$ws.Range("E1").Value = "territory [Link] [Territory]"
$ws.Range("F1").Value = "gender [Link] [Gender]"
$ws.Range("G1").Value = "lead_name [Link] [Lead]"
$ws.Range("H1").Value = "opportunity_name [Link] [Opportunity]"
$ws.Range("I1").Value = "prospect_name [Link] [Prospect]"
$ws.Range("J1").Value = "account_manager [Link] [User]"
$ws.Range("K1").Value = "default_currency [Link] [Currency]"
$ws.Range("L1").Value = "default_bank_account [Link] [Bank Account]"
$ws.Range("M1").Value = "default_price_list [Link] [Price List]"
$ws.Range("O1").Value = "represents_company [Link] [Company]"
$ws.Range("P1").Value = "market_segment [Link] [Market Segment]"
$ws.Range("Q1").Value = "industry [Link] [Industry Type]"
$ws.Range("S1").Value = "language [Link] [Language]"
$ws.Range("U1").Value = "customer_primary_address [Link] [Address]"
$ws.Range("V1").Value = "customer_primary_contact [Link] [Contact]"
$ws.Range("X1").Value = "tax_category [Link] [Tax Category]"
$ws.Range("Y1").Value = "tax_withholding_category [Link] [Tax Withholding Category]"
$ws.Range("Z1").Value = "payment_terms [Link] [Payment Terms Template]"
$ws.Range("AK1").Value = "loyalty_program [Link] [Loyalty Program]"
$ws.Range("AV1").Value = "default_sales_partner [Link] [Sales Partner]"
